# "Generate Report for Archive"
#
# Localization status moved on from the handoff stage: every cell that
# held the shared string "Ready for handoff" now reads "In Translation".
# That string is shared across the Overview sheet (E2:F3) and the two
# per-locale sheets (zh-cn!C2:C3, de-de!C2:C3), so updating each of those
# cells collapses back onto a single shared-string table entry, exactly
# like the source diff shows.
#
# Because the new text is shorter than the old text, the columns that
# were sized to fit it shrink too (Overview!E:F and the "Status" column
# on each locale sheet) - re-apply AutoFit-equivalent sizing there.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Columns shrink to fit the now-shorter "In Translation" text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: "Status" column (C), rows 2-3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = 12.5
